$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "De volta para o futuro"
$ws.Range("D2").Value = "Ficção Científica"
$ws.Range("A3").Value = "Vivendo a vida adoidado"
$ws.Range("D3").Value = "Comédia"
